$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing row 17 (C17/D17 were empty numeric placeholders) ----
$ws.Range("C17").Value = "23:14:51"
$ws.Range("D17").Value = "1.7 Hours"

# ---- Add new row 18 ----
# A18 holds a date-like string ("2026-02-04"); force it to be stored as
# literal text (not auto-converted to a date serial number) by applying a
# text number format before assigning the value.
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "2026-02-04"

$ws.Range("B18").Value = "00:06:07"
$ws.Range("C18").Value = "00:06:10"
$ws.Range("D18").Value = "0 Hours"

# ---- Add new row 19 (totals row) ----
$ws.Range("C19").Value = "Total Duration:"
$ws.Range("D19").Value = "25.5 Hours"

# ---- Match the formatting of the rest of the table (style used by row 17) ----
# Copy the format from an existing row-17 cell and paste only the formatting
# onto the newly added cells so they pick up the same cell style (s="2")
# without disturbing the text values just entered. Row 18 gets all 4
# columns formatted, but row 19 only uses columns C and D.
$ws.Range("A17:D17").Copy() | Out-Null
$ws.Range("A18:D18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C17:D17").Copy() | Out-Null
$ws.Range("C19:D19").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
